$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 8.324917666666666
$ws.Range("H2").Value = 24.974753
$ws.Range("I2").Value = 0.8193616330571973
$ws.Range("J2").Value = 0.8193616330571972
$ws.Range("M2").Value = 3.241087666666667
$ws.Range("N2").Value = 9.723262999999999
$ws.Range("O2").Value = 0.02486257877280725
$ws.Range("P2").Value = 0.02486257877280725
$ws.Range("Q2").Value = 26.98178797544877
$ws.Range("R2").Value = 242.836091779039
$ws.Range("S2").Value = 0.02037144314530056
$ws.Range("T2").Value = 0.02037144314530056
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 8.324917666666666
$ws.Range("H3").Value = 24.974753
$ws.Range("I3").Value = 0.8193616330571973
$ws.Range("J3").Value = 0.8193616330571972
$ws.Range("O3").Value = 0.02096124117795788
$ws.Range("P3").Value = 0.02096124117795788
$ws.Range("Q3").Value = 22.74791244842566
$ws.Range("R3").Value = 204.731212035831
$ws.Range("S3").Value = 0.01717483680247734
$ws.Range("T3").Value = 0.01717483680247734
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 8.324917666666666
$ws.Range("H4").Value = 24.974753
$ws.Range("I4").Value = 0.8193616330571973
$ws.Range("J4").Value = 0.8193616330571972
$ws.Range("M4").Value = 124.3864796666667
$ws.Range("N4").Value = 373.159439
$ws.Range("O4").Value = 0.9541761800492348
$ws.Range("P4").Value = 0.9541761800492349
$ws.Range("Q4").Value = 1035.507202071507
$ws.Range("R4").Value = 9319.564818643568
$ws.Range("S4").Value = 0.7818153531094194
$ws.Range("T4").Value = 0.7818153531094193
$ws.Range("G5").Value = 1.835330666666667
$ws.Range("H5").Value = 5.505992
$ws.Range("I5").Value = 0.1806383669428028
$ws.Range("J5").Value = 0.1806383669428027
$ws.Range("M5").Value = 3.241087666666667
$ws.Range("N5").Value = 9.723262999999999
$ws.Range("O5").Value = 0.02486257877280725
$ws.Range("P5").Value = 0.02486257877280725
$ws.Range("Q5").Value = 5.948467587988445
$ws.Range("R5").Value = 53.536208291896
$ws.Range("S5").Value = 0.004491135627506694
$ws.Range("T5").Value = 0.004491135627506694
$ws.Range("G6").Value = 1.835330666666667
$ws.Range("H6").Value = 5.505992
$ws.Range("I6").Value = 0.1806383669428028
$ws.Range("J6").Value = 0.1806383669428027
$ws.Range("O6").Value = 0.02096124117795788
$ws.Range("P6").Value = 0.02096124117795788
$ws.Range("Q6").Value = 5.015057564642666
$ws.Range("R6").Value = 45.13551808178399
$ws.Range("S6").Value = 0.003786404375480544
$ws.Range("T6").Value = 0.003786404375480543
$ws.Range("G7").Value = 1.835330666666667
$ws.Range("H7").Value = 5.505992
$ws.Range("I7").Value = 0.1806383669428028
$ws.Range("J7").Value = 0.1806383669428027
$ws.Range("M7").Value = 124.3864796666667
$ws.Range("N7").Value = 373.159439
$ws.Range("O7").Value = 0.9541761800492348
$ws.Range("P7").Value = 0.9541761800492349
$ws.Range("Q7").Value = 228.2903206509431
$ws.Range("R7").Value = 2054.612885858488
$ws.Range("S7").Value = 0.1723608269398155
$ws.Range("T7").Value = 0.1723608269398155
